# Auto-generated edit script: updates cached market-price derived values
# in the Anima_Profits tracker sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4800.2173
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 4800.2173
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 14400.6519
$ws.Range("M112").Value = ""
$ws.Range("N112").Value = -16616.6519
$ws.Range("H113").Value = 2767.2632
$ws.Range("I113").Value = 2692.8
$ws.Range("J113").Value = 2850
$ws.Range("K113").Value = 2692.8
$ws.Range("L113").Value = 2850
$ws.Range("M113").Value = 561.1999999999998
$ws.Range("N113").Value = -9358
$ws.Range("H116").Value = 18400.666
$ws.Range("J116").Value = 4000
$ws.Range("L116").Value = 4000
$ws.Range("N116").Value = -10884
$ws.Range("H129").Value = 1895.6389
$ws.Range("I129").Value = 727.7143
$ws.Range("K129").Value = 2183.1429
$ws.Range("M129").Value = 2816.8571
$ws.Range("H132").Value = 2327.739
$ws.Range("I132").Value = 1578.5
$ws.Range("J132").Value = 7322.6665
$ws.Range("K132").Value = 4735.5
$ws.Range("L132").Value = 21967.9995
$ws.Range("M132").Value = -2205.5
$ws.Range("N132").Value = -27027.9995
$ws.Range("H141").Value = 6095
$ws.Range("I141").Value = 3118.75
$ws.Range("J141").Value = 12047.5
$ws.Range("K141").Value = 9356.25
$ws.Range("L141").Value = 36142.5
$ws.Range("M141").Value = -4176.25
$ws.Range("N141").Value = -46502.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1029865.5
$ws.Range("I32").Value = 1088530
$ws.Range("K32").Value = 1088530
$ws.Range("M32").Value = -1088243
$ws.Range("H61").Value = 2972.5652
$ws.Range("I61").Value = 2598.0625
$ws.Range("J61").Value = 3828.5715
$ws.Range("K61").Value = 2598.0625
$ws.Range("L61").Value = 3828.5715
$ws.Range("M61").Value = -2386.0625
$ws.Range("N61").Value = -4252.5715
$ws.Range("H97").Value = 1749.2307
$ws.Range("I97").Value = 1742.5
$ws.Range("J97").Value = 1760
$ws.Range("K97").Value = 1742.5
$ws.Range("L97").Value = 1760
$ws.Range("M97").Value = -1246.5
$ws.Range("N97").Value = -2752
$ws.Range("H122").Value = 35559.133
$ws.Range("I122").Value = 57287.445
$ws.Range("K122").Value = 171862.335
$ws.Range("M122").Value = -169412.335
$ws.Range("H136").Value = 2972.5652
$ws.Range("I136").Value = 2598.0625
$ws.Range("J136").Value = 3828.5715
$ws.Range("K136").Value = 7794.1875
$ws.Range("L136").Value = 11485.7145
$ws.Range("M136").Value = -5244.1875
$ws.Range("N136").Value = -16585.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2749.122
$ws.Range("I134").Value = 2538.2812
$ws.Range("J134").Value = 3498.7778
$ws.Range("K134").Value = 7614.8436
$ws.Range("L134").Value = 10496.3334
$ws.Range("M134").Value = -5079.8436
$ws.Range("N134").Value = -15566.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2387.8845
$ws.Range("I58").Value = 3509.3635
$ws.Range("J58").Value = 1565.4667
$ws.Range("K58").Value = 3509.3635
$ws.Range("L58").Value = 1565.4667
$ws.Range("M58").Value = -3306.3635
$ws.Range("N58").Value = -1971.4667
$ws.Range("H132").Value = 1858.2051
$ws.Range("I132").Value = 1388.2609
$ws.Range("J132").Value = 2533.75
$ws.Range("K132").Value = 4164.7827
$ws.Range("L132").Value = 7601.25
$ws.Range("M132").Value = -1634.7827
$ws.Range("N132").Value = -12661.25
$ws.Range("H134").Value = 6253847
$ws.Range("I134").Value = 10004935
$ws.Range("J134").Value = 2032.6666
$ws.Range("K134").Value = 30014805
$ws.Range("L134").Value = 6097.9998
$ws.Range("M134").Value = -30012270
$ws.Range("N134").Value = -11167.9998
$ws.Range("H136").Value = 2387.8845
$ws.Range("I136").Value = 3509.3635
$ws.Range("J136").Value = 1565.4667
$ws.Range("K136").Value = 10528.0905
$ws.Range("L136").Value = 4696.4001
$ws.Range("M136").Value = -7978.0905
$ws.Range("N136").Value = -9796.400099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 19000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 19000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 57000
$ws.Range("M81").Value = ""
$ws.Range("N81").Value = -59246
$ws.Range("H82").Value = 1000
$ws.Range("I82").Value = 1000
$ws.Range("K82").Value = 3000
$ws.Range("M82").Value = -2594
$ws.Range("H84").Value = 19000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 19000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 171000
$ws.Range("M84").Value = ""
$ws.Range("N84").Value = -182232
$ws.Range("H85").Value = 1000
$ws.Range("I85").Value = 1000
$ws.Range("K85").Value = 3000
$ws.Range("M85").Value = -1596
$ws.Range("H86").Value = 700
$ws.Range("J86").Value = 700
$ws.Range("L86").Value = 2100
$ws.Range("N86").Value = -4472
$ws.Range("H87").Value = 750
$ws.Range("I87").Value = 750
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 2250
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -1002
$ws.Range("N87").Value = ""
$ws.Range("H89").Value = 700
$ws.Range("J89").Value = 700
$ws.Range("L89").Value = 6300
$ws.Range("N89").Value = -18156
$ws.Range("H90").Value = 750
$ws.Range("I90").Value = 750
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 6750
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -510
$ws.Range("N90").Value = ""
$ws.Range("H131").Value = 3662.932
$ws.Range("J131").Value = 4195
$ws.Range("L131").Value = 12585
$ws.Range("N131").Value = -22665
$ws.Range("H139").Value = 4252.0264
$ws.Range("J139").Value = 4954.4644
$ws.Range("L139").Value = 14863.3932
$ws.Range("N139").Value = -25143.3932

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2493.28
$ws.Range("I136").Value = 2091.65
$ws.Range("K136").Value = 6274.950000000001
$ws.Range("M136").Value = -3724.950000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5116.6665
$ws.Range("I136").Value = 5200.7856
$ws.Range("J136").Value = 5026.077
$ws.Range("K136").Value = 15602.3568
$ws.Range("L136").Value = 15078.231
$ws.Range("M136").Value = -13052.3568
$ws.Range("N136").Value = -20178.231
